# Updates cryptos list values (price + 1h volume %) to match the Oct 26
# 2024 GitHub Actions refresh. Rows 42-45 also swap coin order
# (RenderToken/PolygonEcosystemToken and OKB/dogwifhat trade places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is purely numeric-looking ("583.50", "1.00", ...)
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a real number (e.g. "9.01" -> 9.0099999999999998,
# "1.00" -> 1) and the literal formatting from the source data is lost.

$ws.Range("D2").Value = "67.176.29"
$ws.Range("E2").Value = "  -0.98%  "

$ws.Range("D3").Value = "2.470.51"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.50"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.57"
$ws.Range("E6").Value = "  -3.65%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  -2.27%  "

$ws.Range("D9").Value = "2.471.43"
$ws.Range("E9").Value = "  -1.11%  "

$ws.Range("E10").Value = "  -4.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.93"
$ws.Range("E12").Value = "  -3.55%  "

$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.48"
$ws.Range("E15").Value = "  -3.17%  "

$ws.Range("D16").Value = "67.089.06"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("E17").Value = "  -4.75%  "

$ws.Range("D18").Value = "2.487.12"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  -2.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.66"
$ws.Range("E20").Value = "  -4.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.95"
$ws.Range("E21").Value = "  -2.58%  "

$ws.Range("E22").Value = "  -2.75%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.59"
$ws.Range("E24").Value = "  -2.51%  "

$ws.Range("E25").Value = "  -7.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.77"
$ws.Range("E26").Value = "  -7.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.01"
$ws.Range("E27").Value = "  -8.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").Value = "2.594.45"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("E30").Value = "  -6.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "513.28"
$ws.Range("E31").Value = "  -4.18%  "

$ws.Range("E32").Value = "  -5.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.78"
$ws.Range("E33").Value = "  -4.71%  "

$ws.Range("E34").Value = "  -5.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  -7.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.71"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.59"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.47"
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("E40").Value = "  -6.20%  "

$ws.Range("E41").Value = "  -6.57%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.79"
$ws.Range("E42").Value = "  -6.64%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.325"
$ws.Range("E43").Value = "  -6.92%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.67"
$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("E45").Value = "  -7.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "141.52"
$ws.Range("E46").Value = "  -2.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.47"
$ws.Range("E47").Value = "  -5.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.516"
$ws.Range("E48").Value = "  -6.00%  "

$ws.Range("E49").Value = "  -7.43%  "

$ws.Range("E50").Value = "  -6.09%  "

$ws.Range("E51").Value = "  -2.01%  "
